# Update automàtic: dades i banners [2026-02-05 12:00]
# Refreshes the DATA_EXTRACCIO timestamps (column E) for rows 2-36 and
# updates a handful of re-measured weather figures on row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DATA_EXTRACCIO values per row (row number -> timestamp string)
$timestamps = @{
    2  = "2026-02-05 11:58:56"
    3  = "2026-02-05 11:58:59"
    4  = "2026-02-05 11:59:01"
    5  = "2026-02-05 11:59:04"
    6  = "2026-02-05 11:59:07"
    7  = "2026-02-05 11:59:10"
    8  = "2026-02-05 11:59:12"
    9  = "2026-02-05 11:59:15"
    10 = "2026-02-05 11:59:18"
    11 = "2026-02-05 11:59:20"
    12 = "2026-02-05 11:59:23"
    13 = "2026-02-05 11:59:25"
    14 = "2026-02-05 11:59:28"
    15 = "2026-02-05 11:59:31"
    16 = "2026-02-05 11:59:33"
    17 = "2026-02-05 11:59:36"
    18 = "2026-02-05 11:59:39"
    19 = "2026-02-05 11:59:41"
    20 = "2026-02-05 11:59:44"
    21 = "2026-02-05 11:59:46"
    22 = "2026-02-05 11:59:49"
    23 = "2026-02-05 11:59:52"
    24 = "2026-02-05 11:59:54"
    25 = "2026-02-05 11:59:57"
    26 = "2026-02-05 12:00:00"
    27 = "2026-02-05 12:00:03"
    28 = "2026-02-05 12:00:05"
    29 = "2026-02-05 12:00:08"
    30 = "2026-02-05 12:00:10"
    31 = "2026-02-05 12:00:13"
    32 = "2026-02-05 12:00:16"
    33 = "2026-02-05 12:00:18"
    34 = "2026-02-05 12:00:21"
    35 = "2026-02-05 12:00:24"
    36 = "2026-02-05 12:00:27"
}

foreach ($row in $timestamps.Keys) {
    $ws.Range("E$row").Value = $timestamps[$row]
}

# Row 17 (La Seu D'urgell - la Seu D'urgell - Bellestar) also got refreshed
# observation figures alongside its extraction timestamp.
$ws.Range("I17").Value = "2.8 mm"
$ws.Range("J17").Value = "997.4 hPa"
$ws.Range("K17").Value = "0.7 MJ/m2"
$ws.Range("O17").Value = "0.4 °C"
